# Applies the cryptos list update (values for Feb 23 2024 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.208.50"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'2.957.91"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'381.83"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'102.51"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'3.425.18"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'18.04"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'7.42"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "'2.946.94"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "'0.992"
$ws.Range("E17").Value = "  +4.36%  "
$ws.Range("D18").Value = "'51.157.34"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  -6.22%  "
$ws.Range("D20").Value = "'7.13"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("D22").Value = "'0.0₃0954"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'68.48"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'262.08"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'8.41"
$ws.Range("E26").Value = "  +13.33%  "
$ws.Range("D27").Value = "'7.76"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +8.49%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'0.0455"
$ws.Range("E33").Value = "  +5.00%  "
$ws.Range("D34").Value = "'33.95"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.40"
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'16.78"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "'121.60"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'21.29"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'2.009.86"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'0.0346"
$ws.Range("E50").Value = "  +7.17%  "
$ws.Range("D51").Value = "'2.15"
$ws.Range("E51").Value = "  +16.23%  "
